# Apply updated crypto price/volume data (scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.263.59"
$ws.Range("E2").Value = "  -2.81%  "
$ws.Range("D3").Value = "3.680.40"
$ws.Range("E3").Value = "  -3.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'683.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.32%  "
$ws.Range("D6").Value = "'162.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.65%  "
$ws.Range("D7").Value = "3.678.49"
$ws.Range("E7").Value = "  -3.65%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -4.78%  "
$ws.Range("E10").Value = "  -8.54%  "
$ws.Range("D11").Value = "'7.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.26%  "
$ws.Range("E12").Value = "  -3.92%  "
$ws.Range("D13").Value = "'0.0000240"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.66%  "
$ws.Range("D14").Value = "'33.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.72%  "
$ws.Range("D15").Value = "4.301.70"
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Value = "3.678.00"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "69.350.33"
$ws.Range("E17").Value = "  -2.62%  "
$ws.Range("D19").Value = "'16.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.55%  "
$ws.Range("D20").Value = "'6.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.89%  "
$ws.Range("D21").Value = "'482.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.53%  "
$ws.Range("D22").Value = "'9.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.10%  "
$ws.Range("E23").Value = "  -8.25%  "
$ws.Range("D24").Value = "'80.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("D25").Value = "3.826.56"
$ws.Range("E25").Value = "  -3.55%  "
$ws.Range("E26").Value = "  -10.01%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'11.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.49%  "
$ws.Range("D29").Value = "'9.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.83%  "
$ws.Range("E30").Value = "  -10.49%  "
$ws.Range("D31").Value = "'2.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.89%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.71%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'2.08"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.94%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'27.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.29%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.167"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "3.651.12"
$ws.Range("E38").Value = "  -7.62%  "
$ws.Range("D39").Value = "'6.31"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.67%  "
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").Value = "'0.0934"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.08%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'0.949"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.34%  "
$ws.Range("D45").Value = "'163.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.46%  "
$ws.Range("D46").Value = "'48.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.84%  "
$ws.Range("E47").Value = "  -13.80%  "
$ws.Range("D48").Value = "'29.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("D49").Value = "'0.000286"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.64%  "
$ws.Range("D50").Value = "'1.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("E51").Value = "  -3.36%  "
